$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.178.96"
$ws.Range("E2").Value = "  -2.19%  "

$ws.Range("D3").Value = "3.131.45"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.18%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "3.123.60"
$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("E9").Value = "  -2.68%  "

$ws.Range("E10").Value = "  -3.41%  "

$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.68%  "

$ws.Range("D15").Value = "3.644.73"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").Value = "63.134.79"
$ws.Range("E17").Value = "  -2.17%  "

$ws.Range("D18").Value = "3.130.63"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("E19").Value = "  -2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("E21").Value = "  -4.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.92%  "

$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("E25").Value = "  -4.09%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("E28").Value = "  -3.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.41%  "

$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("E33").Value = "  -5.81%  "

$ws.Range("E34").Value = "  -4.40%  "

$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("E36").Value = "  -2.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").Value = "0.0₃0712"
$ws.Range("E38").Value = "  -4.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "424.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.79%  "

$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("E42").Value = "  -10.11%  "

$ws.Range("D43").Value = "2.894.25"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.85%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("E49").Value = "  -5.96%  "

$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
